# add odoxa poll (11/11)
# Appends 3 new poll rows (id=30, odoxa, online, partially-unsure, 11/11/2021,
# week 10) to the bottom of the Sheet1 data table, one row per right-wing
# candidate scenario (Bertrand / Pecresse / Barnier).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=id B=year C=week D=month E=day F=firm G=collectmode H=unsure
#          I=n J=c_poutou K=c_arthaud L=c_melenchon M=c_roussel N=c_montebourg
#          O=c_jadot P=c_hidalgo Q=c_macron R=c_pecresse S=c_barnier
#          T=c_bertrand U=c_lassalle V=c_daignant W=c_lepen X=c_zemmour

# Row 102 - Bertrand scenario
$ws.Range("A102").Value = 30
$ws.Range("B102").Value = 2021
$ws.Range("C102").Value = 10
$ws.Range("D102").Value = 11
$ws.Range("E102").Value = 7
$ws.Range("F102").Value = "odoxa"
$ws.Range("G102").Value = "online"
$ws.Range("H102").Value = "partially"
$ws.Range("I102").Value = 1917
$ws.Range("J102").Value = 2
$ws.Range("K102").Value = 1
$ws.Range("L102").Value = 8.5
$ws.Range("M102").Value = 2
$ws.Range("N102").Value = 2
$ws.Range("O102").Value = 6.5
$ws.Range("P102").Value = 4.5
$ws.Range("Q102").Value = 25
$ws.Range("T102").Value = 12
$ws.Range("U102").Value = 1.5
$ws.Range("V102").Value = 3
$ws.Range("W102").Value = 18
$ws.Range("X102").Value = 14

# Row 103 - Pecresse scenario
$ws.Range("A103").Value = 30
$ws.Range("B103").Value = 2021
$ws.Range("C103").Value = 10
$ws.Range("D103").Value = 11
$ws.Range("E103").Value = 7
$ws.Range("F103").Value = "odoxa"
$ws.Range("G103").Value = "online"
$ws.Range("H103").Value = "partially"
$ws.Range("I103").Value = 1946
$ws.Range("J103").Value = 2.5
$ws.Range("K103").Value = 1
$ws.Range("L103").Value = 8.5
$ws.Range("M103").Value = 2
$ws.Range("N103").Value = 2
$ws.Range("O103").Value = 6.5
$ws.Range("P103").Value = 5
$ws.Range("Q103").Value = 25
$ws.Range("R103").Value = 9
$ws.Range("U103").Value = 1.5
$ws.Range("V103").Value = 4
$ws.Range("W103").Value = 18.5
$ws.Range("X103").Value = 14.5

# Row 104 - Barnier scenario
$ws.Range("A104").Value = 30
$ws.Range("B104").Value = 2021
$ws.Range("C104").Value = 10
$ws.Range("D104").Value = 11
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = "odoxa"
$ws.Range("G104").Value = "online"
$ws.Range("H104").Value = "partially"
$ws.Range("I104").Value = 1917
$ws.Range("J104").Value = 2
$ws.Range("K104").Value = 1
$ws.Range("L104").Value = 9
$ws.Range("M104").Value = 2
$ws.Range("N104").Value = 2.5
$ws.Range("O104").Value = 7
$ws.Range("P104").Value = 5
$ws.Range("Q104").Value = 25
$ws.Range("S104").Value = 9
$ws.Range("U104").Value = 1
$ws.Range("V104").Value = 3
$ws.Range("W104").Value = 19
$ws.Range("X104").Value = 14.5

# Leave the selection on the last entered cell, matching the authored view.
$ws.Range("U104").Select()
